$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NATMI cell-cell communication metrics (Cxcl13-Cxcr3) with
# recomputed values from the new TPM expression data.

# Row 2
$ws.Range("I2").Value = 0.5162107379131895
$ws.Range("J2").Value = 0.5162107379131895
$ws.Range("O2").Value = 0.02773017886769741
$ws.Range("P2").Value = 0.02773017886769741
$ws.Range("S2").Value = 0.01431461609575881
$ws.Range("T2").Value = 0.01431461609575881

# Row 3
$ws.Range("I3").Value = 0.5162107379131895
$ws.Range("J3").Value = 0.5162107379131895
$ws.Range("M3").Value = 1.821156333333333
$ws.Range("N3").Value = 5.463469
$ws.Range("O3").Value = 0.9722698211323025
$ws.Range("P3").Value = 0.9722698211323026
$ws.Range("Q3").Value = 7.369245969413777
$ws.Range("R3").Value = 66.323213724724
$ws.Range("S3").Value = 0.5018961218174307
$ws.Range("T3").Value = 0.5018961218174307

# Row 4
$ws.Range("G4").Value = 0.155986
$ws.Range("H4").Value = 0.467958
$ws.Range("I4").Value = 0.01989925565426652
$ws.Range("J4").Value = 0.01989925565426652
$ws.Range("O4").Value = 0.02773017886769741
$ws.Range("P4").Value = 0.02773017886769741
$ws.Range("Q4").Value = 0.008102120821333332
$ws.Range("R4").Value = 0.072919087392
$ws.Range("S4").Value = 0.0005518099186268496
$ws.Range("T4").Value = 0.0005518099186268496

# Row 5
$ws.Range("G5").Value = 0.155986
$ws.Range("H5").Value = 0.467958
$ws.Range("I5").Value = 0.01989925565426652
$ws.Range("J5").Value = 0.01989925565426652
$ws.Range("M5").Value = 1.821156333333333
$ws.Range("N5").Value = 5.463469
$ws.Range("O5").Value = 0.9722698211323025
$ws.Range("P5").Value = 0.9722698211323026
$ws.Range("Q5").Value = 0.2840748918113333
$ws.Range("R5").Value = 2.556674026302
$ws.Range("S5").Value = 0.01934744573563967
$ws.Range("T5").Value = 0.01934744573563967

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.636334333333334
$ws.Range("H6").Value = 10.909003
$ws.Range("I6").Value = 0.463890006432544
$ws.Range("J6").Value = 0.463890006432544
$ws.Range("O6").Value = 0.02773017886769741
$ws.Range("P6").Value = 0.02773017886769741
$ws.Range("Q6").Value = 0.1888760537191111
$ws.Range("R6").Value = 1.699884483472
$ws.Range("S6").Value = 0.01286375285331175
$ws.Range("T6").Value = 0.01286375285331175

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.636334333333334
$ws.Range("H7").Value = 10.909003
$ws.Range("I7").Value = 0.463890006432544
$ws.Range("J7").Value = 0.463890006432544
$ws.Range("M7").Value = 1.821156333333333
$ws.Range("N7").Value = 5.463469
$ws.Range("O7").Value = 0.9722698211323025
$ws.Range("P7").Value = 0.9722698211323026
$ws.Range("Q7").Value = 6.622333301267445
$ws.Range("R7").Value = 59.600999711407
$ws.Range("S7").Value = 0.4510262535792323
$ws.Range("T7").Value = 0.4510262535792323
